$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:I70 and J2:J70
$iValues = @(9,9,9,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,7,7,8,9,6,9,6,7,9,9,9,8,8,7,9,9,9,9,7,9,8,9,9,9,9,8,9,8,9,9,8,7,9,8,9,8,9,8,9,7,6,7,8,7,6,3)
$jValues = @(9,10,10,9,9,9,9,9,9,10,10,9,9,9,9,9,9,9,9,10,9,9,9,7,7,9,9,6,10,7,7,9,9,9,9,8,7,9,9,9,9,7,9,8,9,9,9,9,8,10,8,9,9,8,7,9,8,9,8,10,8,9,7,7,7,8,7,6,3)

for ($r = 2; $r -le 70; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
